$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values remain text (matches original inlineStr formatting),
# since many look like numbers (e.g. 0.9996) and Excel would otherwise convert them.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.067.06'
$ws.Range('E2').Value = '  +2.60%  '
$ws.Range('D3').Value = '1.814.84'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  -0.48%  '
$ws.Range('D5').Value = '337.94'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').Value = '0.9958'
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('D7').Value = '0.3925'
$ws.Range('E7').Value = '  +3.15%  '
$ws.Range('D8').Value = '0.3486'
$ws.Range('E8').Value = '  +0.96%  '
$ws.Range('D9').Value = '48.44'
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('D11').Value = '0.07591'
$ws.Range('E11').Value = '  +0.85%  '
$ws.Range('D12').Value = '0.9977'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D13').Value = '22.19'
$ws.Range('E13').Value = '  +0.31%  '
$ws.Range('D14').Value = '6.543'
$ws.Range('E14').Value = '  +0.67%  '
$ws.Range('D15').Value = '1.814.10'
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').Value = '7.198'
$ws.Range('E16').Value = '  +1.57%  '
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('D18').Value = '0.06705'
$ws.Range('E18').Value = '  +0.63%  '
$ws.Range('D19').Value = '85.39'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').Value = '0.9973'
$ws.Range('D21').Value = '17.85'
$ws.Range('E21').Value = '  +2.66%  '
$ws.Range('D22').Value = '6.574'
$ws.Range('E22').Value = '  +0.52%  '
$ws.Range('D23').Value = '28.083.43'
$ws.Range('E23').Value = '  +2.65%  '
$ws.Range('D24').Value = '12.85'
$ws.Range('E24').Value = '  +2.36%  '
$ws.Range('D25').Value = '2.407'
$ws.Range('E25').Value = '  -0.92%  '
$ws.Range('E26').Value = '  +0.92%  '
$ws.Range('D27').Value = '2.572'
$ws.Range('E27').Value = '  -0.52%  '
$ws.Range('D28').Value = '21.37'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('D29').Value = '154.84'
$ws.Range('E29').Value = '  +1.85%  '
$ws.Range('D30').Value = '2.019.74'
$ws.Range('E30').Value = '  +1.14%  '
$ws.Range('D31').Value = '135.87'
$ws.Range('E31').Value = '  +1.15%  '
$ws.Range('D32').Value = '4.042'
$ws.Range('E32').Value = '  -0.55%  '
$ws.Range('D33').Value = '6.149'
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('D34').Value = '0.08848'
$ws.Range('E34').Value = '  +1.44%  '
$ws.Range('D35').Value = '13.30'
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('D36').Value = '5.527'
$ws.Range('E36').Value = '  +1.21%  '
$ws.Range('D37').Value = '0.6951'
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('D38').Value = '0.02428'
$ws.Range('E38').Value = '  +3.86%  '
$ws.Range('D39').Value = '0.06553'
$ws.Range('E39').Value = '  +2.67%  '
$ws.Range('D40').Value = '1.616'
$ws.Range('E40').Value = '  -4.35%  '
$ws.Range('E41').Value = '  +1.08%  '
$ws.Range('D42').Value = '1.265'
$ws.Range('E42').Value = '  -0.88%  '
$ws.Range('D43').Value = '8.535'
$ws.Range('E43').Value = '  -4.10%  '
$ws.Range('D44').Value = '14.64'
$ws.Range('E44').Value = '  +1.26%  '
$ws.Range('D45').Value = '0.6507'
$ws.Range('E45').Value = '  +0.37%  '
$ws.Range('D46').Value = '0.9959'
$ws.Range('E46').Value = '  -0.47%  '
$ws.Range('D47').Value = '3.873'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('D48').Value = '2.170'
$ws.Range('E48').Value = '  +2.01%  '
$ws.Range('D49').Value = '132.62'
$ws.Range('E49').Value = '  +1.68%  '
$ws.Range('D50').Value = '0.07220'
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('D51').Value = '80.38'
$ws.Range('E51').Value = '  +1.20%  '
